# Generate Report for Handoff
# b.md has now been prepared and sent out for handoff in both locales, so
# update its status / handoff file / handoff datetime on each sheet.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: row 3 corresponds to b.md ---
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-03-24 09:36:37"

# --- zh-cn sheet: row 3 corresponds to b.md ---
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("D3").Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("E3").Value = "2016-03-24 09:36:32"

# --- de-de sheet: row 3 corresponds to b.md ---
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("D3").Hyperlinks.Item(1).TextToDisplay = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("E3").Value = "2016-03-24 09:36:37"
